$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing hours entries (row 3, 4, 6) ---
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 6
$ws.Range("B6").Value = 2

# --- Apply the existing date style (style index used by A2:A19, numFmtId 14)
#     to the new date cells by copying format only, so the new rows reuse
#     the same cellXf instead of creating a new one. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A20:A29").PasteSpecial(-4122) | Out-Null

# --- New timesheet rows, written in the same order the original author
#     entered them (controls shared-string allocation order) ---

# Row 21/22: "Testing in work week" (shared twice)
$ws.Range("A21").Value = 41303
$ws.Range("B21").Value = 6
$ws.Range("C21").Value = "Testing in work week"

$ws.Range("A22").Value = 41304
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = "Testing in work week"

# Row 23
$ws.Range("A23").Value = 41305
$ws.Range("B23").Value = 14
$ws.Range("C23").Value = "Testing in work week/also completing test documentation"

# Row 24
$ws.Range("A24").Value = 41306
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "Final testing"

# Row 20
$ws.Range("A20").Value = 41302
$ws.Range("B20").Value = 6
$ws.Range("C20").Value = "Testing in work week/also imrpoving on design documentation"

# Row 25 (reuses existing "Group Meeting" shared string)
$ws.Range("A25").Value = 41312
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "Group Meeting"

# Row 26
$ws.Range("A26").Value = 41317
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Completing my personal report"

# Row 27
$ws.Range("A27").Value = 41317
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = "Completing the change in control forms"

# Row 28
$ws.Range("A28").Value = 41317
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "Complete the timesheet documentation"

# Row 29 (note: K29/M29 already hold "Matt Whitmore" / "Group 17" and must survive)
$ws.Range("A29").Value = 41317
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = "Comleting my timesheet"

# Row 31: grand total
$ws.Range("A31").Value = "Total"
$ws.Range("B31").Value = 72

# --- View state: selection moves to B25, with the viewport scrolled so
#     row 9 is at the top (best-effort - ScrollRow is honoured by the
#     object model even where the xlsx writer may not round-trip it). ---
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("B25").Select() | Out-Null
